# Auto-applied numeric corrections to market-price / profit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, driven by a
# scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38 (marker G38=4599)
$ws.Range("H38").Value = 870
$ws.Range("I38").Value = 870
$ws.Range("K38").Value = 2610
$ws.Range("M38").Value = -2238

# Row 40 (marker G40=5505)
$ws.Range("H40").Value = 2830
$ws.Range("I40").Value = 2660.3333
$ws.Range("J40").Value = 2999.6667
$ws.Range("K40").Value = 2660.3333
$ws.Range("L40").Value = 2999.6667
$ws.Range("M40").Value = -2485.3333
$ws.Range("N40").Value = -3349.6667

# Row 51 (marker G51=5486)
$ws.Range("H51").Value = 5593.8
$ws.Range("J51").Value = 5742.25
$ws.Range("L51").Value = 5742.25
$ws.Range("N51").Value = -6710.25

# Row 112 (marker G112=27960)
$ws.Range("H112").Value = 5322.0835
$ws.Range("J112").Value = 6146.5
$ws.Range("L112").Value = 18439.5
$ws.Range("N112").Value = -20655.5

# Row 129 (marker G129=36115)
$ws.Range("H129").Value = 904.64386
$ws.Range("I129").Value = 1203.8
$ws.Range("J129").Value = 882.64703
$ws.Range("K129").Value = 3611.4
$ws.Range("L129").Value = 2647.94109
$ws.Range("M129").Value = 1388.6
$ws.Range("N129").Value = -12647.94109

# Row 137 (marker G137=44013)
$ws.Range("H137").Value = 1951.8667
$ws.Range("J137").Value = 2131.5
$ws.Range("L137").Value = 6394.5
$ws.Range("N137").Value = -11494.5

# Row 138 (marker G138=44169)
$ws.Range("H138").Value = 2611.3103
$ws.Range("I138").Value = 2736.0435
$ws.Range("J138").Value = 2133.1667
$ws.Range("K138").Value = 8208.130500000001
$ws.Range("L138").Value = 6399.500100000001
$ws.Range("M138").Value = -3068.130500000001
$ws.Range("N138").Value = -16679.5001

$ws = $wb.Worksheets.Item("ARM")
# Row 19 (marker G19=3550)
$ws.Range("H19").Value = 10004
$ws.Range("I19").Value = 10004
$ws.Range("K19").Value = 10004
$ws.Range("M19").Value = -9775

# Row 30 (marker G30=2712)
$ws.Range("H30").Value = 7777
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 7777
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 7777
$ws.Range("N30").Value = -8077
$ws.Range("M30").ClearContents()

# Row 32 (marker G32=44147)
$ws.Range("H32").Value = 4925.477
$ws.Range("I32").Value = 3594.0667
$ws.Range("K32").Value = 3594.0667
$ws.Range("M32").Value = -3307.0667

# Row 74 (marker G74=44000)
$ws.Range("H74").Value = 481.33334
$ws.Range("I74").Value = 481.33334
$ws.Range("K74").Value = 481.33334
$ws.Range("M74").Value = 392.66666

# Row 77 (marker G77=44000)
$ws.Range("H77").Value = 481.33334
$ws.Range("I77").Value = 481.33334
$ws.Range("K77").Value = 2406.6667
$ws.Range("M77").Value = 1961.3333

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (marker G22=5092)
$ws.Range("H22").Value = 463.66666
$ws.Range("I22").Value = 356.6
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 356.6
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -183.6
$ws.Range("N22").Value = -1345

# Row 95 (marker G95=18194)
$ws.Range("H95").Value = 71897.336
$ws.Range("J95").Value = 71897.336
$ws.Range("L95").Value = 71897.336
$ws.Range("N95").Value = -77389.336

$ws = $wb.Worksheets.Item("CRP")
# Row 6 (marker G6=2219)
$ws.Range("H6").Value = 6250
$ws.Range("I6").Value = 2500
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 2500
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = -2387
$ws.Range("N6").Value = -10226

# Row 31 (marker G31=44023)
$ws.Range("H31").Value = 2393.9688
$ws.Range("I31").Value = 1152.9584
$ws.Range("K31").Value = 1152.9584
$ws.Range("M31").Value = -857.9584

# Row 34 (marker G34=44023)
$ws.Range("H34").Value = 2393.9688
$ws.Range("I34").Value = 1152.9584
$ws.Range("K34").Value = 1152.9584
$ws.Range("M34").Value = -950.9584

$ws = $wb.Worksheets.Item("CUL")
# Row 23 (marker G23=4858)
$ws.Range("H23").Value = 148.125
$ws.Range("I23").Value = 49.5
$ws.Range("J23").Value = 181
$ws.Range("K23").Value = 148.5
$ws.Range("L23").Value = 543
$ws.Range("M23").Value = 86.5
$ws.Range("N23").Value = -1013

# Row 131 (marker G131=36060)
$ws.Range("H131").Value = 8351.947
$ws.Range("I131").Value = 534.3333
$ws.Range("J131").Value = 8884.966
$ws.Range("K131").Value = 1602.9999
$ws.Range("L131").Value = 26654.898
$ws.Range("M131").Value = 3437.0001
$ws.Range("N131").Value = -36734.898

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (marker G80=12521)
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3000
$ws.Range("N80").Value = -4996
$ws.Range("M80").ClearContents()

# Row 83 (marker G83=12521)
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 15000
$ws.Range("N83").Value = -24984
$ws.Range("M83").ClearContents()

# Row 97 (marker G97=19940)
$ws.Range("H97").Value = 2759.8
$ws.Range("I97").Value = 2350
$ws.Range("J97").Value = 3033
$ws.Range("K97").Value = 2350
$ws.Range("L97").Value = 3033
$ws.Range("M97").Value = -1854
$ws.Range("N97").Value = -4025

# Row 102 (marker G102=36169)
$ws.Range("H102").Value = 2227.1035
$ws.Range("I102").Value = 2293.9375
$ws.Range("J102").Value = 2144.8462
$ws.Range("K102").Value = 2293.9375
$ws.Range("L102").Value = 2144.8462
$ws.Range("M102").Value = -671.9375
$ws.Range("N102").Value = -5388.8462

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (marker G22=5277)
$ws.Range("H22").Value = 1981
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 27 (marker G27=5277)
$ws.Range("H27").Value = 1981
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# Row 68 (marker G68=12563)
$ws.Range("H68").Value = 3773.2222
$ws.Range("I68").Value = 3568.4285
$ws.Range("K68").Value = 3568.4285
$ws.Range("M68").Value = -2819.4285

# Row 71 (marker G71=12563)
$ws.Range("H71").Value = 3773.2222
$ws.Range("I71").Value = 3568.4285
$ws.Range("K71").Value = 17842.1425
$ws.Range("M71").Value = -14098.1425

# Row 82 (marker G82=12565)
$ws.Range("H82").Value = 4384.75
$ws.Range("I82").Value = 1499.5
$ws.Range("J82").Value = 5346.5
$ws.Range("K82").Value = 1499.5
$ws.Range("L82").Value = 5346.5
$ws.Range("M82").Value = -1138.5
$ws.Range("N82").Value = -6068.5

# Row 85 (marker G85=12565)
$ws.Range("H85").Value = 4384.75
$ws.Range("I85").Value = 1499.5
$ws.Range("J85").Value = 5346.5
$ws.Range("K85").Value = 1499.5
$ws.Range("L85").Value = 5346.5
$ws.Range("M85").Value = -251.5
$ws.Range("N85").Value = -7842.5

# Row 98 (marker G98=18379)
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990

# Row 100 (marker G100=19995)
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (marker G107=27746)
$ws.Range("H107").Value = 608.8182
$ws.Range("I107").Value = 410.77777
$ws.Range("K107").Value = 1232.33331
$ws.Range("M107").Value = 687.66669

$wb.Save()